# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (col I) and DialogAct (col J) values for specific rows in Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 9).Value = "b"
$ws.Cells.Item(4, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(10, 9).Value = "sv"
$ws.Cells.Item(10, 10).Value = "Statement-opinion"
$ws.Cells.Item(14, 9).Value = "b"
$ws.Cells.Item(14, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(15, 9).Value = "sd"
$ws.Cells.Item(15, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(22, 9).Value = "aa"
$ws.Cells.Item(22, 10).Value = "Agree/Accept"
$ws.Cells.Item(26, 9).Value = "sv"
$ws.Cells.Item(26, 10).Value = "Statement-opinion"
$ws.Cells.Item(28, 9).Value = "sd"
$ws.Cells.Item(28, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(37, 9).Value = "sv"
$ws.Cells.Item(37, 10).Value = "Statement-opinion"
$ws.Cells.Item(39, 9).Value = "aa"
$ws.Cells.Item(39, 10).Value = "Agree/Accept"
$ws.Cells.Item(40, 9).Value = "sv"
$ws.Cells.Item(40, 10).Value = "Statement-opinion"
$ws.Cells.Item(44, 9).Value = "aa"
$ws.Cells.Item(44, 10).Value = "Agree/Accept"
$ws.Cells.Item(57, 9).Value = "ba"
$ws.Cells.Item(57, 10).Value = "Appreciation"
$ws.Cells.Item(65, 9).Value = "sv"
$ws.Cells.Item(65, 10).Value = "Statement-opinion"
$ws.Cells.Item(66, 9).Value = "sd"
$ws.Cells.Item(66, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(68, 9).Value = "sd"
$ws.Cells.Item(68, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(87, 9).Value = "sv"
$ws.Cells.Item(87, 10).Value = "Statement-opinion"
$ws.Cells.Item(90, 9).Value = "sd"
$ws.Cells.Item(90, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(100, 9).Value = "b"
$ws.Cells.Item(100, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(114, 9).Value = "aa"
$ws.Cells.Item(114, 10).Value = "Agree/Accept"
$ws.Cells.Item(115, 9).Value = "aa"
$ws.Cells.Item(115, 10).Value = "Agree/Accept"
$ws.Cells.Item(131, 9).Value = "sd"
$ws.Cells.Item(131, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(143, 9).Value = "sd"
$ws.Cells.Item(143, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(146, 9).Value = "sd"
$ws.Cells.Item(146, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(152, 9).Value = "sv"
$ws.Cells.Item(152, 10).Value = "Statement-opinion"
$ws.Cells.Item(158, 9).Value = "aa"
$ws.Cells.Item(158, 10).Value = "Agree/Accept"
$ws.Cells.Item(159, 9).Value = "aa"
$ws.Cells.Item(159, 10).Value = "Agree/Accept"
$ws.Cells.Item(167, 9).Value = "aa"
$ws.Cells.Item(167, 10).Value = "Agree/Accept"
$ws.Cells.Item(176, 9).Value = "sv"
$ws.Cells.Item(176, 10).Value = "Statement-opinion"
$ws.Cells.Item(191, 9).Value = "sd"
$ws.Cells.Item(191, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(194, 9).Value = "sv"
$ws.Cells.Item(194, 10).Value = "Statement-opinion"
$ws.Cells.Item(196, 9).Value = "sv"
$ws.Cells.Item(196, 10).Value = "Statement-opinion"
$ws.Cells.Item(201, 9).Value = "qy"
$ws.Cells.Item(201, 10).Value = "Yes-No-Question"
$ws.Cells.Item(206, 9).Value = "%"
$ws.Cells.Item(206, 10).Value = "Uninterpretable"
$ws.Cells.Item(227, 9).Value = "sv"
$ws.Cells.Item(227, 10).Value = "Statement-opinion"
$ws.Cells.Item(248, 9).Value = "aa"
$ws.Cells.Item(248, 10).Value = "Agree/Accept"
$ws.Cells.Item(249, 9).Value = "aa"
$ws.Cells.Item(249, 10).Value = "Agree/Accept"
$ws.Cells.Item(256, 9).Value = "aa"
$ws.Cells.Item(256, 10).Value = "Agree/Accept"
$ws.Cells.Item(272, 9).Value = "sd"
$ws.Cells.Item(272, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(275, 9).Value = "sd"
$ws.Cells.Item(275, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(277, 9).Value = "sd"
$ws.Cells.Item(277, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(282, 9).Value = "aa"
$ws.Cells.Item(282, 10).Value = "Agree/Accept"
$ws.Cells.Item(287, 9).Value = "sv"
$ws.Cells.Item(287, 10).Value = "Statement-opinion"
$ws.Cells.Item(289, 9).Value = "sv"
$ws.Cells.Item(289, 10).Value = "Statement-opinion"
$ws.Cells.Item(291, 9).Value = "sd"
$ws.Cells.Item(291, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(302, 9).Value = "sd"
$ws.Cells.Item(302, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(327, 9).Value = "ba"
$ws.Cells.Item(327, 10).Value = "Appreciation"
$ws.Cells.Item(371, 9).Value = "sv"
$ws.Cells.Item(371, 10).Value = "Statement-opinion"
$ws.Cells.Item(372, 9).Value = "sv"
$ws.Cells.Item(372, 10).Value = "Statement-opinion"
